$d = $word.ActiveDocument

# Update the title/date line (first paragraph).
$d.Content.Find.Execute("2024-12-14 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-12-15 Sunday", 2) | Out-Null

# Update the division problems inside the table. Addressed by explicit
# row/cell index (rather than Find/Replace) because several of the
# original values (e.g. "25÷2=", "21÷3=") repeat, and Find would be
# ambiguous about which occurrence to touch.
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("67÷6=", "93÷7=", "48÷8=", "14÷3=", "10÷8=")
    5  = @("98÷5=", "96÷3=", "37÷2=", "25÷8=", "53÷7=")
    9  = @("63÷9=", "19÷9=", "10÷3=", "21÷3=", "81÷7=")
    13 = @("46÷8=", "94÷9=", "56÷9=", "94÷4=", "99÷3=")
    17 = @("23÷5=", "58÷5=", "93÷9=", "82÷5=", "47÷5=")
}

foreach ($rowIndex in $updates.Keys) {
    $row = $t.Rows.Item($rowIndex)
    $values = $updates[$rowIndex]
    for ($c = 1; $c -le $values.Length; $c++) {
        $row.Cells.Item($c).Range.Text = $values[$c - 1]
    }
}
